# 12/9/2017 MAMATHA AND CHITRA CHICK IN
$d = $word.ActiveDocument

# --- Change 1: merge the split "Mon Sep 10" / " 12:18:54 PDT 2017" runs
# into a single run by re-finding & replacing the paragraph text (Word
# COM Find/Replace rewrites the found range as a single run when the
# replacement text shares one formatting).
$d.Content.Find.Execute(
    "Mon Sep 10 12:18:54 PDT 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "Mon Sep 10 12:18:54 PDT 2017", 2) | Out-Null

# --- Change 2: append a new "Tue Sep 11" purchase record after the
# final "Amount balance" paragraph (the last populated paragraph in the
# document), pushing it ahead of the trailing blank paragraphs.
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$courier = "<w:rFonts w:ascii='Courier New' w:hAnsi='Courier New' w:cs='Courier New'/>"
$courierB = "<w:rFonts w:ascii='Courier New' w:hAnsi='Courier New' w:cs='Courier New'/><w:b/>"
$courierRed = "<w:rFonts w:ascii='Courier New' w:hAnsi='Courier New' w:cs='Courier New'/><w:color w:val='FF0000'/>"

$xml = ""

# blank bold paragraph (spacer before the new date header)
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courierB</w:rPr></w:pPr></w:p>"

# date header: "Tue Sep 11 12:31:52 PDT 2017" (two runs, like the source records)
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courier</w:rPr></w:pPr>" +
        "<w:r><w:rPr>$courier</w:rPr><w:t>Tue Sep 11</w:t></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:t xml:space='preserve'> 12:31:52 PDT 2017</w:t></w:r>" +
        "</w:p>"

# dashed separator line
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courier</w:rPr></w:pPr>" +
        "<w:r><w:rPr>$courier</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r>" +
        "</w:p>"

# Person Name - YASHODHA
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courier</w:rPr></w:pPr>" +
        "<w:r><w:rPr>$courier</w:rPr><w:t>Person Name</w:t></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:tab/></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:tab/></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:tab/></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:tab/><w:t>- YASHODHA</w:t></w:r>" +
        "</w:p>"

# Amount Received - 1924 (red)
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courierRed</w:rPr></w:pPr>" +
        "<w:r><w:rPr>$courierRed</w:rPr><w:t>Amount Received</w:t></w:r>" +
        "<w:r><w:rPr>$courierRed</w:rPr><w:tab/></w:r>" +
        "<w:r><w:rPr>$courierRed</w:rPr><w:tab/></w:r>" +
        "<w:r><w:rPr>$courierRed</w:rPr><w:tab/><w:t>- 1924</w:t></w:r>" +
        "</w:p>"

# Amount Received mode - CASH AND CLEARD
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courier</w:rPr></w:pPr>" +
        "<w:r><w:rPr>$courier</w:rPr><w:t>Amount Received mode</w:t></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:tab/></w:r>" +
        "<w:r><w:rPr>$courier</w:rPr><w:tab/><w:t>- CASH AND CLEARD</w:t></w:r>" +
        "</w:p>"

# trailing blank bold paragraph
$xml += "<w:p $w><w:pPr><w:pStyle w:val='PlainText'/><w:rPr>$courierB</w:rPr></w:pPr></w:p>"

# Locate the last "Amount balance" paragraph (the final populated record)
# and insert the new block right after it, before the existing trailing
# blank paragraphs.
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Amount balance*") {
        $target = $p
        break
    }
}

$ins = $d.Range($target.Range.End, $target.Range.End)
$ins.InsertXML($xml) | Out-Null
